# RMI updates through 12/2
# Replace several formula-driven "percent of fuel demand change" cells on the
# PoFDCtAE sheet with static value 1 (100%), matching the author's manual override.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PoFDCtAE")

$cells = @("C3", "D4", "I9", "J10", "K11", "L12", "M13", "N14", "S19", "T20")
foreach ($cellRef in $cells) {
    $ws.Range($cellRef).Value = 1
}

$excel.Calculate()
